$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for rows 2-9 (weekly update: new row added, data re-sorted by date desc)
$data = @(
    @{ Row=2; D=44699; L="Primera"; M=100; N=20000; O=22000; P=21000; Q="$/caja 18 kilos";        S=1167 },
    @{ Row=3; D=44699; L="Segunda"; M=50;  N=18000; O=18000; P=18000; Q="$/caja 18 kilos";        S=1000 },
    @{ Row=4; D=44819; L="Primera"; M=100; N=25000; O=26000; P=25500; Q="$/caja 18 kilos granel"; S=1417 },
    @{ Row=5; D=44516; L="Primera"; M=100; N=33000; O=34000; P=33500; Q="$/caja 18 kilos";        S=1861 },
    @{ Row=6; D=44687; L="Primera"; M=100; N=18000; O=19000; P=18500; Q="$/caja 18 kilos";        S=1028 },
    @{ Row=7; D=44316; L="Segunda"; M=50;  N=20000; O=20000; P=20000; Q="$/caja 18 kilos";        S=1111 },
    @{ Row=8; D=44280; L="Primera"; M=100; N=14000; O=15000; P=14500; Q="$/caja 18 kilos";        S=806  },
    @{ Row=9; D=44280; L="Segunda"; M=50;  N=12000; O=12000; P=12000; Q="$/caja 18 kilos";        S=667  }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 4).Value  = $row.D   # D: Fecha
    $ws.Cells.Item($r, 12).Value = $row.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $row.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $row.N   # N: Precio mínimo
    $ws.Cells.Item($r, 15).Value = $row.O   # O: Precio máximo
    $ws.Cells.Item($r, 16).Value = $row.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row.Q   # Q: Unidad de comercialización
    $ws.Cells.Item($r, 19).Value = $row.S   # S: Precio $/Kg
}
